# Cloning onto new PC.
#
# The last slide ("A bit of history on compilers" / References slide) had
# its body-copy recap line edited: "A bit of history on compilers" became
# "A bit of history of compilers", with the edited word ("history of ")
# now living in its own run (distinct formatting history from the rest of
# the sentence, which was left untouched).

$p = $ppt.ActivePresentation

# The slide is the very last one in the deck (sldId 314 / creationId
# {FC4A93C0-0CAD-4BA5-9658-F523BC4C732F} on the body placeholder).
$s = $p.Slides.Item($p.Slides.Count)

# Shape 1 = Title ("A bit of history on compilers" -- left untouched).
# Shape 2 = "Content Placeholder 2", which holds the recap paragraph.
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 2 is "A bit of history on compilers<br>https://...".
$para = $tr.Paragraphs(2, 1)
$run = $para.Runs(1)

# Replace the "history on " substring (chars 10-20, 1-based) in the
# original run with "history of ", splitting the run into three: the
# untouched lead-in, the edited middle word, and the untouched tail.
$mid = $run.Characters(10, 11)
$mid.Text = "history of "
